# Avoid insertion of records if it already has the new items in the Queue to Process
#
# Replace the sample Transactions data with a smaller set (15 rows instead of
# 25), including a couple of repeated "already queued" item names (abc / acc)
# used to exercise the "don't re-add existing queue items" logic, clear the
# now-unused trailing rows, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data table (rows 1-16). Column A mixes numeric transaction
# ids with a few repeated text ids ("abc"/"acc") that already exist in the
# queue; column D is the "UseCashCount" mode.
$data = @(
    @("CashIn", "OnUsCheck", "NotOnUsCheck", "UseCashCount"),
    @(110, 510, 1, "Use Both"),
    @(120, 520, 2, "Use Piece Count"),
    @(130, 530, 3, "Use Piece Count"),
    @(140, 540, 4, "Use Amount"),
    @("abc", 550, 5, "Use Both"),
    @(160, 560, 6, "Use Both"),
    @(170, 570, 7, "Use Both"),
    @("acc", 580, 8, "Use Amount"),
    @(190, 590, 9, "Use Amount"),
    @(200, 600, 10, "Use Piece Count"),
    @(210, 610, 11, "Use Piece Count"),
    @(220, 620, 12, "Use Amount"),
    @("acc", 630, 13, "Use Both"),
    @(240, 640, 14, "Use Piece Count"),
    @(250, 650, 15, "Use Both")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Rows 17-32 previously held extra sample rows (up to row 26) plus blank
# filler rows (27-32); now they all become fully empty (only the formatted
# A/B placeholder cells remain, matching rows 33+).
$null = $ws.Range("A17:D26").ClearContents()
$null = $ws.Range("C17:D32").Clear()

# Update the active selection to match the saved workbook state.
$null = $ws.Range("C16").Select()
